# Fix Rstudio -> RStudio
#
# The workbook has a single worksheet with a schedule; column D holds the
# "topic" for each row. Row 6 (cell D6) contains the topic text
# "Introduction to R, Rstudio, and the tidyverse" which needs the product
# name corrected to "RStudio" (capital S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Introduction to R, RStudio, and the tidyverse"

# Update the active selection to match the authored commit (D7 instead of D12)
$ws.Range("D7").Select()
